# report-checklist.xlsx update
# - Column K ("Tipo Documento") gets filled in with "Tipo Documento non gestito"
#   for all the test-case rows that previously had it blank/varied.
# - The AutoFilter (filtering on column J = "SI") is removed, which also
#   reveals every row that had been hidden by the filter.
# - The active selection / scroll position is moved to A26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "Tipo Documento non gestito"

# Rows whose "K" cell (Tipo Documento) must be set to the unified text.
# These are exactly the rows that were hidden by the old AutoFilter.
$kRowRanges = @(
    @(8, 21),
    @(23, 25),
    @(27, 87),
    @(105, 109)
)

foreach ($range in $kRowRanges) {
    for ($r = $range[0]; $r -le $range[1]; $r++) {
        $ws.Range("K$r").Value = $newText
    }
}

# Remove the AutoFilter entirely (also unhides the filtered-out rows).
$ws.AutoFilterMode = $false

# Move the selection / view to A26.
$ws.Range("A26").Select()
